$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - ALGASON MASSAGE CREAM 40 GM: current balance changed 1:0 -> 2:0
$ws.Range("H7").Value = "2:0"

# Row 13 - item changed from "IVY PRONT  SYRUP" to "IVYPRONT 0.84GM SYRUP 120 ML"
# with updated balance / order-limit / price / selling-price
$ws.Range("C13").Value = "IVYPRONT 0.84GM SYRUP 120 ML"
$ws.Range("H13").Value = "1:0"
$ws.Range("L13").NumberFormat = "@"
$ws.Range("L13").Value = "1"
$ws.Range("N13").Value = "45.00"
$ws.Range("P13").NumberFormat = "@"
$ws.Range("P13").Value = "45.0000"

# Row 17 - REPAION-N GEL 50 GM: current balance changed 1:0 -> 2:0
$ws.Range("H17").Value = "2:0"

# Grand total reflects the price change on row 13 (57.00 -> 45.00 = -12.00)
$ws.Range("P22").Value = 623.98

# Footer timestamp updated
$ws.Range("A23").Value = "Thursday, 9 October, 2025 12:48 PM"
